$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.187.98"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "'3.170.12"
$ws.Range("E3").Value = "  +3.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'573.41"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").Value = "'150.82"
$ws.Range("E6").Value = "  +5.97%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.166.61"
$ws.Range("E8").Value = "  +3.84%  "
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +5.17%  "
$ws.Range("D11").Value = "'6.19"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  +6.13%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  +18.46%  "
$ws.Range("D14").Value = "'38.25"
$ws.Range("E14").Value = "  +8.86%  "
$ws.Range("D15").Value = "'3.686.14"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "'65.230.77"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "'3.182.13"
$ws.Range("E18").Value = "  +7.22%  "
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "'512.34"
$ws.Range("E20").Value = "  +7.43%  "
$ws.Range("D21").Value = "'14.96"
$ws.Range("E21").Value = "  +7.19%  "
$ws.Range("D22").Value = "'16.05"
$ws.Range("E22").Value = "  +12.65%  "
$ws.Range("D23").Value = "'0.736"
$ws.Range("E23").Value = "  +8.45%  "
$ws.Range("D24").Value = "'7.86"
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("D25").Value = "'84.94"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  +15.54%  "
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("E29").Value = "  +9.28%  "
$ws.Range("D30").Value = "'28.07"
$ws.Range("E30").Value = "  +7.11%  "
$ws.Range("E31").Value = "  +15.48%  "
$ws.Range("E32").Value = "  +8.07%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +12.47%  "
$ws.Range("D35").Value = "'6.68"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("D36").Value = "'55.60"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "'477.26"
$ws.Range("E37").Value = "  +7.86%  "
$ws.Range("E38").Value = "  +9.54%  "
$ws.Range("D39").Value = "'3.09"
$ws.Range("E39").Value = "  +9.07%  "
$ws.Range("E40").Value = "  +4.33%  "
$ws.Range("D41").Value = "'3.131.50"
$ws.Range("E41").Value = "  +5.12%  "
$ws.Range("D42").Value = "'8.64"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("D43").Value = "'0.121"
$ws.Range("E43").Value = "  +7.39%  "
$ws.Range("E45").Value = "  +11.76%  "
$ws.Range("D46").Value = "'29.25"
$ws.Range("E46").Value = "  +5.92%  "
$ws.Range("E47").Value = "  +13.98%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("E50").Value = "  +12.66%  "
$ws.Range("D51").Value = "'123.60"
$ws.Range("E51").Value = "  +5.41%  "
